$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new backtest entry (date 2017-07-29 / serial 42945)
# Copy the formatting (date style) from the row above (A5) into A6, then set its value.
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 42945

$ws.Range("B6").Value = 0.00065137

$ws.Range("P6").Formula = "=SUM(B6:L6)"

$ws.Range("R6").Value = 2677

$ws.Range("S6").Formula = "=P6*R6"

# Update the active selection to match the saved view state
$ws.Range("C6").Select()
